$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Tipo de revision:" row (row 4) - its label row, shifting rows 5-34 up by one
$ws.Rows.Item(4).Delete()

# After the first delete, the "Inspector(es):" row, originally row 7, is now row 6
$ws.Rows.Item(6).Delete()

# Update the title cell text
$ws.Range("B2").Value = "Reporte"

# Update the selection to match the target state
$ws.Range("B7").Select()
